# Se incorpora la url del repositorio en la presentación.
#
# Adds a new final slide (same "Título y contenido" layout/style as the
# rest of the deck) whose body placeholder contains the project's GitHub
# repository URL as a clickable hyperlink.

$p = $ppt.ActivePresentation

# Duplicate one of the deck's plain "Título y contenido" slides so the
# new one inherits the normal slide formatting/structure (layout,
# color-map override, empty shape/body properties, etc.) instead of the
# bare placeholders a brand-new slide would get. Slide 6 is a clean,
# minimal instance of that layout (no leftover size/autofit overrides,
# no bullet/indent formatting on its first body paragraph).
# Duplicate() inserts right after its source, so move the copy to the
# end of the deck afterwards.
$templateSlide = $p.Slides.Item(6)
$dupRange = $templateSlide.Duplicate()
$slide = $dupRange.Item(1)
$slide.MoveTo($p.Slides.Count)

# --- Title placeholder -------------------------------------------------
$title = $slide.Shapes.Item(1)
$title.Name = "Título 4"
$title.TextFrame.TextRange.Text = "Guía Completa"
$title.TextFrame.TextRange.LanguageID = "es-CL"

# --- Body / content placeholder ----------------------------------------
$body = $slide.Shapes.Item(2)
$body.Name = "Marcador de contenido 5"

$prefix = "Repositorio: "
$url = "https://github.com/jmcandia/python_para_desarrolladores"

$content = $body.TextFrame.TextRange
$content.Text = $prefix + $url
$content.LanguageID = "es-CL"

$linkRange = $content.Characters($prefix.Length + 1, $url.Length)
$linkRange.ActionSettings.Item(1).Hyperlink.Address = $url
